$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.75, 4933110.244521, 46529.238143, 4886581.006379, 130759.439386, 1401001.109276, 22013.885332, 1378987.223945, 139076.682658),
    @(0.8, 5057225.336225, 52431.195682, 5004794.140543, 158493.584141, 1536912.016744, 30256.719986, 1506655.296759, 167416.276471),
    @(0.85, 5887332.721518, 53041.027801, 5834291.693717, 161614.921388, 2530383.478741, 32404.702314, 2497978.776427, 170678.598832)
)

$startRow = 18
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
